$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(1)
Write-Host $sh.Name
try {
    $sh.Fill.ForeColor.ObjectThemeColor = 5  # msoThemeColorAccent1 or similar
    Write-Host "set ok"
} catch { Write-Host "ERR $_" }
